# Auto update Excel log
# Appends the latest door-sensor event to the Proximity log and the
# matching camera-capture event to the Camera log.

$wb = $excel.ActiveWorkbook

# --- Proximity sheet: new row 10 (EXIT event) ---------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
$rowProximity = 10

# Column A holds a date-shaped string ("2026-02-01"). Writing it straight
# through .Value would get auto-converted to a date serial by Excel's
# input parser, so force Text formatting for the write, then clear the
# formatting back to the sheet's default so no stray style is left behind.
$wsProximity.Range("A$rowProximity").NumberFormat = "@"
$wsProximity.Range("A$rowProximity").Value = "2026-02-01"
$wsProximity.Range("A$rowProximity").ClearFormats()

$wsProximity.Range("B$rowProximity").Value = "16:45:11"
$wsProximity.Range("C$rowProximity").Value = "16:00"
$wsProximity.Range("D$rowProximity").Value = "Living Room Main Door"
$wsProximity.Range("E$rowProximity").Value = "EXIT"
$wsProximity.Range("F$rowProximity").Value = "User EXITED Living Room Main Door"

# --- Camera sheet: new row 9 (Image Captured event) ----------------------
$wsCamera = $wb.Worksheets.Item("Camera")
$rowCamera = 9

$wsCamera.Range("A$rowCamera").NumberFormat = "@"
$wsCamera.Range("A$rowCamera").Value = "2026-02-01"
$wsCamera.Range("A$rowCamera").ClearFormats()

$wsCamera.Range("B$rowCamera").Value = "16:45:11"
$wsCamera.Range("C$rowCamera").Value = "16:00"
$wsCamera.Range("D$rowCamera").Value = "Living Room Main Door"
$wsCamera.Range("E$rowCamera").Value = "Image Captured"
$wsCamera.Range("F$rowCamera").Value = "Active"
